# Actualización automática 2025-08-19 12:20:09
$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M13").Value = 5034.91

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 5034.91
$ws2.Range("F23").Value = 6556.13

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 6315.17
$ws3.Range("E16").Value = 32461.3
$ws3.Range("F16").Value = 0.1628608792909721
$ws3.Range("D19").Value = 6556.13
$ws3.Range("E19").Value = 52832.09762291769
$ws3.Range("F19").Value = 0.1103944377937626

$ws3.Columns.Item(6).ColumnWidth = 23.166666666666668
